$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $origStyle = $ws.Range($addr).Style
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = $origStyle
}

Set-TextValue $ws "D2" "44.619.04"
Set-TextValue $ws "E2" "  +0.96%  "
Set-TextValue $ws "D3" "2.245.57"
Set-TextValue $ws "E3" "  +0.16%  "
Set-TextValue $ws "E4" "  +0.28%  "
Set-TextValue $ws "D5" "306.36"
Set-TextValue $ws "E5" "  -0.06%  "
Set-TextValue $ws "D6" "94.67"
Set-TextValue $ws "E6" "  +0.09%  "
Set-TextValue $ws "D7" "0.569"
Set-TextValue $ws "E7" "  -0.34%  "
Set-TextValue $ws "E8" "  +0.05%  "
Set-TextValue $ws "D9" "0.516"
Set-TextValue $ws "E9" "  -1.72%  "
Set-TextValue $ws "D10" "34.80"
Set-TextValue $ws "E10" "  +0.37%  "
Set-TextValue $ws "D11" "0.0799"
Set-TextValue $ws "E11" "  -1.38%  "
Set-TextValue $ws "E12" "  +0.28%  "
Set-TextValue $ws "E13" "  +0.07%  "
Set-TextValue $ws "D14" "2.588.36"
Set-TextValue $ws "E14" "  +0.02%  "
Set-TextValue $ws "D15" "2.241.13"
Set-TextValue $ws "E15" "  -3.69%  "
Set-TextValue $ws "E16" "  +0.20%  "
Set-TextValue $ws "D17" "13.56"
Set-TextValue $ws "E17" "  +0.28%  "
Set-TextValue $ws "D18" "44.398.67"
Set-TextValue $ws "E18" "  +0.99%  "
Set-TextValue $ws "D19" "0.0₃0935"
Set-TextValue $ws "E19" "  -2.83%  "
Set-TextValue $ws "E20" "  -3.06%  "
Set-TextValue $ws "D21" "11.75"
Set-TextValue $ws "E21" "  -3.07%  "
Set-TextValue $ws "D22" "65.31"
Set-TextValue $ws "E22" "  -0.34%  "
Set-TextValue $ws "D23" "237.68"
Set-TextValue $ws "E23" "  -0.02%  "
Set-TextValue $ws "E24" "  -0.07%  "
Set-TextValue $ws "E25" "  -1.35%  "
Set-TextValue $ws "E26" "  -0.05%  "
Set-TextValue $ws "E27" "  +4.24%  "
Set-TextValue $ws "E28" "  -1.63%  "
Set-TextValue $ws "D29" "36.91"
Set-TextValue $ws "E29" "  -3.39%  "
Set-TextValue $ws "D30" "19.99"
Set-TextValue $ws "E30" "  -0.15%  "
Set-TextValue $ws "D31" "5.85"
Set-TextValue $ws "E31" "  +0.21%  "
Set-TextValue $ws "D32" "147.48"
Set-TextValue $ws "E32" "  -3.71%  "
Set-TextValue $ws "D33" "0.0781"
Set-TextValue $ws "E33" "  -1.71%  "
Set-TextValue $ws "E34" "  +0.08%  "
Set-TextValue $ws "E35" "  +1.08%  "
Set-TextValue $ws "E36" "  +1.22%  "
Set-TextValue $ws "E37" "  -1.51%  "
Set-TextValue $ws "E38" "  +5.36%  "
Set-TextValue $ws "D39" "15.20"
Set-TextValue $ws "E39" "  +5.39%  "
Set-TextValue $ws "E40" "  -4.94%  "
Set-TextValue $ws "D41" "3.76"
Set-TextValue $ws "E41" "  -1.44%  "
Set-TextValue $ws "E42" "  +0.15%  "
Set-TextValue $ws "E43" "  +0.05%  "
Set-TextValue $ws "D44" "1.810.03"
Set-TextValue $ws "E44" "  +3.78%  "
Set-TextValue $ws "D45" "1.76"
Set-TextValue $ws "E45" "  +11.52%  "
Set-TextValue $ws "D46" "81.97"
Set-TextValue $ws "E46" "  -0.67%  "
Set-TextValue $ws "E47" "  -1.74%  "
Set-TextValue $ws "D48" "98.19"
Set-TextValue $ws "E48" "  -1.54%  "
Set-TextValue $ws "D49" "68.87"
Set-TextValue $ws "E49" "  +2.71%  "
Set-TextValue $ws "D50" "4.81"
Set-TextValue $ws "E50" "  -2.45%  "
Set-TextValue $ws "D51" "53.99"
Set-TextValue $ws "E51" "  -0.70%  "
